# ------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a new row into the "总计" (summary) sheet for 2022-Q4.
# 2. Insert a brand-new worksheet named "2022-Q4" right after "总计"
#    and fill it with the per-fund holdings data for that quarter.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------- 1. Update the "总计" summary sheet ----------------------
$summary = $wb.Worksheets.Item(1)

# Make room for the new quarter as the new row 2 (right under the header)
$summary.Rows.Item(2).Insert()

# Copy formatting (borders/font/alignment) from the row that used to be
# row 2 (now shifted to row 3) down onto the freshly inserted row so the
# new row looks consistent with the rest of the table.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.09

# Renumber the index column (A) for the rows that shifted down
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------- 2. Add the new "2022-Q4" worksheet ----------------------
# Keep a handle on an already-existing, fully-styled quarter sheet so we
# can copy its look & feel (borders/bold header/centered index column)
# onto the brand new sheet instead of inventing a new style from blank
# cells.
$styleSource = $wb.Worksheets.Item("2022-Q2")

$afterSheet = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $afterSheet)
$q4.Name = "2022-Q4"

# Replicate sheetPr / outline settings and overall layout by copying the
# whole used range formatting from the sibling sheet first.
$styleSource.Range("A1:H6").Copy()
$q4.Range("A1:H6").PasteSpecial(-4122)
$q4.Range("A1:H6").ClearContents()

# Header row (row 1), starting at column B just like the sibling sheets
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

function Set-FundRow {
    param($ws, $row, $idx, $code, $fundName, $scale, $stockPos, $posRatio, $marketValue, $posRank)

    $ws.Range("A$row").Value = $idx

    # Columns B..G hold text that merely looks numeric (fund codes with
    # leading zeros, percentages, decimal amounts); force a text number
    # format first so Excel does not silently coerce them to doubles.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $code

    $ws.Range("C$row").Value = $fundName

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $scale

    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $stockPos

    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = $posRatio

    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = $marketValue

    $ws.Range("H$row").Value = $posRank
}

Set-FundRow $q4 2 0 "005457" "景顺长城量化小盘股票"          "5.08" "94.36" "1.34" "0.0681" 10
Set-FundRow $q4 3 1 "015496" "景顺中证1000指数增强C"         "0.86" "92.30" "1.28" "0.0110" 9
Set-FundRow $q4 4 2 "015495" "景顺中证1000指数增强A"         "0.67" "92.30" "1.28" "0.0086" 9
Set-FundRow $q4 5 3 "000926" "中信建投睿信灵活配置混合A"     "0.10" "83.25" "1.03" "0.0010" 9
Set-FundRow $q4 6 4 "004676" "中信建投睿信灵活配置混合C"     "0.03" "83.25" "1.03" "0.0003" 9

Write-Host "Sheet count:"
Write-Host $wb.Worksheets.Count
